$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep text formatting so numeric-looking strings
# (e.g. "318.35", "0.630") are not auto-converted to numbers by Excel.
$targetCells = @('D2','E2','D3','E3','E4','D5','E5','D6','E6','D7','E7','E8','E9','D10','E10','D11','E11','D12','E12','E13','E14','D15','E15','D16','E16','D17','E17','D18','E18','D19','E19','E20','D21','E21','D22','E22','E23','D24','E24','E25','E26','E27','E28','D29','E29','D30','E30','B31','C31','D31','E31','B32','C32','D32','E32','D33','E33','E34','D35','E35','E36','E37','E38','D39','E39','E40','D41','E41','D42','E42','D43','E43','E44','E45','D46','E46','D47','E47','D48','E48','E49','D50','E50','D51','E51')
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.750.44'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '2.312.16'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '318.35'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = '104.65'
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').Value = '40.05'
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('D11').Value = '0.0907'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '8.54'
$ws.Range('E12').Value = '  +3.67%  '
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').Value = '15.42'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').Value = '2.662.76'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '2.325.32'
$ws.Range('E17').Value = '  +1.64%  '
$ws.Range('D18').Value = '42.699.03'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').Value = '7.64'
$ws.Range('E19').Value = '  +2.62%  '
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('D21').Value = '13.47'
$ws.Range('E21').Value = '  +34.42%  '
$ws.Range('D22').Value = '74.02'
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').Value = '269.11'
$ws.Range('E24').Value = '  -5.02%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('E27').Value = '  +1.21%  '
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('D29').Value = '22.70'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').Value = '37.91'
$ws.Range('E30').Value = '  +5.68%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '6.25'
$ws.Range('E31').Value = '  +7.38%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '165.62'
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('D33').Value = '0.0892'
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').Value = '2.59'
$ws.Range('E35').Value = '  -8.75%  '
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('E37').Value = '  +2.15%  '
$ws.Range('E38').Value = '  +1.71%  '
$ws.Range('D39').Value = '3.71'
$ws.Range('E39').Value = '  +1.50%  '
$ws.Range('E40').Value = '  -5.98%  '
$ws.Range('D41').Value = '1.64'
$ws.Range('E41').Value = '  +12.91%  '
$ws.Range('D42').Value = '97.84'
$ws.Range('E42').Value = '  -2.59%  '
$ws.Range('D43').Value = '70.37'
$ws.Range('E43').Value = '  +1.39%  '
$ws.Range('E44').Value = '  +1.10%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '12.39'
$ws.Range('E46').Value = '  +4.30%  '
$ws.Range('D47').Value = '116.52'
$ws.Range('E47').Value = '  +1.99%  '
$ws.Range('D48').Value = '82.35'
$ws.Range('E48').Value = '  +7.54%  '
$ws.Range('E49').Value = '  -0.36%  '
$ws.Range('D50').Value = '5.29'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Value = '1.624.87'
$ws.Range('E51').Value = '  +5.19%  '
